$d = $word.ActiveDocument

$replacements = @(
    @{ old = "2025-02-05 Wednesday"; new = "2025-02-06 Thursday" },
    @{ old = "29×70=2030";  new = "13×21=273"  },
    @{ old = "80×98=7840";  new = "72×70=5040" },
    @{ old = "25×20=500";   new = "86×35=3010" },
    @{ old = "17×87=1479";  new = "24×39=936"  },
    @{ old = "32×88=2816";  new = "24×95=2280" },
    @{ old = "12×26=312";   new = "88×57=5016" },
    @{ old = "71×84=5964";  new = "79×78=6162" },
    @{ old = "44×30=1320";  new = "81×97=7857" },
    @{ old = "84×31=2604";  new = "98×76=7448" },
    @{ old = "20×27=540";   new = "33×72=2376" },
    @{ old = "69×79=5451";  new = "88×62=5456" },
    @{ old = "45×96=4320";  new = "16×30=480"  },
    @{ old = "18×21=378";   new = "39×75=2925" },
    @{ old = "96×57=5472";  new = "62×27=1674" },
    @{ old = "85×37=3145";  new = "90×31=2790" },
    @{ old = "62×52=3224";  new = "39×68=2652" },
    @{ old = "58×72=4176";  new = "41×29=1189" },
    @{ old = "98×45=4410";  new = "11×56=616"  },
    @{ old = "45×85=3825";  new = "28×27=756"  },
    @{ old = "73×46=3358";  new = "22×83=1826" },
    @{ old = "60×65=3900";  new = "62×86=5332" },
    @{ old = "97×25=2425";  new = "11×15=165"  },
    @{ old = "84×36=3024";  new = "70×20=1400" },
    @{ old = "21×50=1050";  new = "91×15=1365" },
    @{ old = "46×72=3312";  new = "83×40=3320" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
